$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segregation of QA and STG product codes: rotate the "current" product
# code shown for each configured row (B2:B7) to a freshly generated one.
$ws.Range("B2").Value = "prodtxga"
$ws.Range("B3").Value = "prodWyTq"
$ws.Range("B4").Value = "prodMtZR"
$ws.Range("B5").Value = "prodSinN"
$ws.Range("B6").Value = "prodNXCU"
$ws.Range("B7").Value = "prodrdtA"
